$wb = $excel.ActiveWorkbook

# Map: worksheet name -> { row -> newValue } for column F ("想去人数")
$updates = @{
    "展览" = @{
        5  = 15
        7  = 1628
        9  = 20
        10 = 1438
        12 = 34
        13 = 364
        15 = 184
        16 = 6
        19 = 259
        20 = 142
        22 = 199
    }
    "全部类型" = @{
        5  = 15
        7  = 1628
        10 = 20
        11 = 1438
        13 = 34
        14 = 364
        16 = 184
        17 = 6
        20 = 259
        21 = 142
        23 = 199
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
